$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 832.6667
$ws.Cells.Item(29, 9).Value = 832.6667
$ws.Cells.Item(29, 11).Value = 2498.0001
$ws.Cells.Item(29, 13).Value = -2217.0001

$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 13).ClearContents()

$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).ClearContents()

$ws.Cells.Item(46, 8).Value = 1795
$ws.Cells.Item(46, 9).Value = 1100
$ws.Cells.Item(46, 10).Value = 2490
$ws.Cells.Item(46, 11).Value = 3300
$ws.Cells.Item(46, 12).Value = 7470
$ws.Cells.Item(46, 13).Value = -3181
$ws.Cells.Item(46, 14).Value = -7708

$ws.Cells.Item(60, 8).Value = 1795
$ws.Cells.Item(60, 9).Value = 1100
$ws.Cells.Item(60, 10).Value = 2490
$ws.Cells.Item(60, 11).Value = 3300
$ws.Cells.Item(60, 12).Value = 7470
$ws.Cells.Item(60, 13).Value = -2816
$ws.Cells.Item(60, 14).Value = -8438

$ws.Cells.Item(116, 8).Value = 340183.88
$ws.Cells.Item(116, 9).Value = 1114301.1
$ws.Cells.Item(116, 11).Value = 1114301.1
$ws.Cells.Item(116, 13).Value = -1110859.1

$ws.Cells.Item(132, 8).Value = 6287.087
$ws.Cells.Item(132, 9).Value = 10912.909
$ws.Cells.Item(132, 10).Value = 2046.75
$ws.Cells.Item(132, 11).Value = 32738.727
$ws.Cells.Item(132, 12).Value = 6140.25
$ws.Cells.Item(132, 13).Value = -30208.727
$ws.Cells.Item(132, 14).Value = -11200.25

$ws.Cells.Item(138, 8).Value = 4478.37
$ws.Cells.Item(138, 9).Value = 748.96155
$ws.Cells.Item(138, 10).Value = 5788.7026
$ws.Cells.Item(138, 11).Value = 2246.88465
$ws.Cells.Item(138, 12).Value = 17366.1078
$ws.Cells.Item(138, 13).Value = 2893.11535
$ws.Cells.Item(138, 14).Value = -27646.1078

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 5798.3335
$ws.Cells.Item(3, 10).Value = 10660
$ws.Cells.Item(3, 12).Value = 10660
$ws.Cells.Item(3, 14).Value = -10890

$ws.Cells.Item(32, 8).Value = 4528.281
$ws.Cells.Item(32, 9).Value = 3903.7551
$ws.Cells.Item(32, 11).Value = 3903.7551
$ws.Cells.Item(32, 13).Value = -3616.7551

$ws.Cells.Item(61, 8).Value = 1602.375
$ws.Cells.Item(61, 9).Value = 1649.8182
$ws.Cells.Item(61, 10).Value = 1498
$ws.Cells.Item(61, 11).Value = 1649.8182
$ws.Cells.Item(61, 12).Value = 1498
$ws.Cells.Item(61, 13).Value = -1437.8182
$ws.Cells.Item(61, 14).Value = -1922

$ws.Cells.Item(74, 8).Value = 6690.75
$ws.Cells.Item(74, 9).Value = 10974.857
$ws.Cells.Item(74, 10).Value = 3358.6667
$ws.Cells.Item(74, 11).Value = 10974.857
$ws.Cells.Item(74, 12).Value = 3358.6667
$ws.Cells.Item(74, 13).Value = -10100.857
$ws.Cells.Item(74, 14).Value = -5106.6667

$ws.Cells.Item(77, 8).Value = 6690.75
$ws.Cells.Item(77, 9).Value = 10974.857
$ws.Cells.Item(77, 10).Value = 3358.6667
$ws.Cells.Item(77, 11).Value = 54874.285
$ws.Cells.Item(77, 12).Value = 16793.3335
$ws.Cells.Item(77, 13).Value = -50506.285
$ws.Cells.Item(77, 14).Value = -25529.3335

$ws.Cells.Item(132, 8).Value = 3232.5908
$ws.Cells.Item(132, 9).Value = 2026.5454
$ws.Cells.Item(132, 10).Value = 4438.636
$ws.Cells.Item(132, 11).Value = 6079.6362
$ws.Cells.Item(132, 12).Value = 13315.908
$ws.Cells.Item(132, 13).Value = -3549.6362
$ws.Cells.Item(132, 14).Value = -18375.908

$ws.Cells.Item(136, 8).Value = 1602.375
$ws.Cells.Item(136, 9).Value = 1649.8182
$ws.Cells.Item(136, 10).Value = 1498
$ws.Cells.Item(136, 11).Value = 4949.4546
$ws.Cells.Item(136, 12).Value = 4494
$ws.Cells.Item(136, 13).Value = -2399.4546
$ws.Cells.Item(136, 14).Value = -9594

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2366.6487
$ws.Cells.Item(134, 9).Value = 1617.7778
$ws.Cells.Item(134, 11).Value = 4853.3334
$ws.Cells.Item(134, 13).Value = -2318.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3963.4119
$ws.Cells.Item(31, 9).Value = 1237.4
$ws.Cells.Item(31, 11).Value = 1237.4
$ws.Cells.Item(31, 13).Value = -942.4000000000001

$ws.Cells.Item(34, 8).Value = 3963.4119
$ws.Cells.Item(34, 9).Value = 1237.4
$ws.Cells.Item(34, 11).Value = 1237.4
$ws.Cells.Item(34, 13).Value = -1035.4

$ws.Cells.Item(58, 8).Value = 1918.8644
$ws.Cells.Item(58, 9).Value = 1725.537
$ws.Cells.Item(58, 11).Value = 1725.537
$ws.Cells.Item(58, 13).Value = -1522.537

$ws.Cells.Item(132, 8).Value = 1996.7805
$ws.Cells.Item(132, 9).Value = 1427.1613
$ws.Cells.Item(132, 11).Value = 4281.4839
$ws.Cells.Item(132, 13).Value = -1751.4839

$ws.Cells.Item(134, 8).Value = 9642.071
$ws.Cells.Item(134, 9).Value = 12443.223
$ws.Cells.Item(134, 11).Value = 37329.669
$ws.Cells.Item(134, 13).Value = -34794.669

$ws.Cells.Item(136, 8).Value = 1918.8644
$ws.Cells.Item(136, 9).Value = 1725.537
$ws.Cells.Item(136, 11).Value = 5176.611
$ws.Cells.Item(136, 13).Value = -2626.611

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1207.7142
$ws.Cells.Item(5, 9).Value = 367.9091
$ws.Cells.Item(5, 10).Value = 2131.5
$ws.Cells.Item(5, 11).Value = 1103.7273
$ws.Cells.Item(5, 12).Value = 6394.5
$ws.Cells.Item(5, 13).Value = -991.7273
$ws.Cells.Item(5, 14).Value = -6618.5

$ws.Cells.Item(122, 8).Value = 3143.4
$ws.Cells.Item(122, 9).Value = 1212
$ws.Cells.Item(122, 10).Value = 3715.6667
$ws.Cells.Item(122, 11).Value = 10908
$ws.Cells.Item(122, 12).Value = 33441.0003
$ws.Cells.Item(122, 13).Value = -8458
$ws.Cells.Item(122, 14).Value = -38341.0003

$ws.Cells.Item(135, 8).Value = 1207.7142
$ws.Cells.Item(135, 9).Value = 367.9091
$ws.Cells.Item(135, 10).Value = 2131.5
$ws.Cells.Item(135, 11).Value = 3311.1819
$ws.Cells.Item(135, 12).Value = 19183.5
$ws.Cells.Item(135, 13).Value = -776.1819
$ws.Cells.Item(135, 14).Value = -24253.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 7128.8
$ws.Cells.Item(132, 9).Value = 5666
$ws.Cells.Item(132, 10).Value = 7494.5
$ws.Cells.Item(132, 11).Value = 16998
$ws.Cells.Item(132, 12).Value = 22483.5
$ws.Cells.Item(132, 13).Value = -14468
$ws.Cells.Item(132, 14).Value = -27543.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1500.8572
$ws.Cells.Item(46, 9).Value = 3464
$ws.Cells.Item(46, 10).Value = 1173.6666
$ws.Cells.Item(46, 11).Value = 3464
$ws.Cells.Item(46, 12).Value = 1173.6666
$ws.Cells.Item(46, 13).Value = -3276
$ws.Cells.Item(46, 14).Value = -1549.6666

$ws.Cells.Item(54, 8).Value = 33814.5
$ws.Cells.Item(54, 10).Value = 33814.5
$ws.Cells.Item(54, 12).Value = 33814.5
$ws.Cells.Item(54, 14).Value = -35102.5

$ws.Cells.Item(68, 8).Value = 690.9299999999999
$ws.Cells.Item(68, 9).Value = 690.9299999999999
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 690.9299999999999
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 58.07000000000005
$ws.Cells.Item(68, 14).ClearContents()

$ws.Cells.Item(71, 8).Value = 690.9299999999999
$ws.Cells.Item(71, 9).Value = 690.9299999999999
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 3454.65
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 289.3500000000004
$ws.Cells.Item(71, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 15875588
$ws.Cells.Item(132, 9).Value = 1694.6666
$ws.Cells.Item(132, 10).Value = 55560320
$ws.Cells.Item(132, 11).Value = 5083.9998
$ws.Cells.Item(132, 12).Value = 166680960
$ws.Cells.Item(132, 13).Value = -2553.9998
$ws.Cells.Item(132, 14).Value = -166686020
